$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their original text formatting (avoid numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.675.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.598.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.27"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.601.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.665.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.672"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.295.64"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.843"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.788"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.735.46"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.896"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.09"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.13%  "
